# Insert two new data rows (242-243) into the "Fruta, Vega Monumental
# Concepcion - Pera" sheet, pushing the existing rows 242.. down by two.
# Excel's native Insert() on a multi-row range copies formatting (incl. the
# date number format on column D) from the row(s) above, which matches the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("242:243").Insert()

# --- New row 242: Abate Fettel / Primera, 2021-11-11 ---
$r = $ws.Cells.Item(242, 1)
$r.Value2 = 11
$ws.Cells.Item(242, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(242, 3).Value2 = "Bíobío"
$ws.Cells.Item(242, 4).Value2 = 44511
$ws.Cells.Item(242, 5).Value2 = 8
$ws.Cells.Item(242, 6).Value2 = "Fruta"
$ws.Cells.Item(242, 7).Value2 = 100104
$ws.Cells.Item(242, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(242, 9).Value2 = 100104005
$ws.Cells.Item(242, 10).Value2 = "Pera"
$ws.Cells.Item(242, 11).Value2 = "Abate Fettel"
$ws.Cells.Item(242, 12).Value2 = "Primera"
$ws.Cells.Item(242, 13).Value2 = 250
$ws.Cells.Item(242, 14).Value2 = 7500
$ws.Cells.Item(242, 15).Value2 = 8000
$ws.Cells.Item(242, 16).Value2 = 7700
$ws.Cells.Item(242, 17).Value2 = "`$/caja 17 kilos empedrada"
$ws.Cells.Item(242, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(242, 19).Value2 = 453
$ws.Cells.Item(242, 20).Value2 = 17

# --- New row 243: Packham's Triumph / Primera, 2021-11-11 ---
$ws.Cells.Item(243, 1).Value2 = 11
$ws.Cells.Item(243, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(243, 3).Value2 = "Bíobío"
$ws.Cells.Item(243, 4).Value2 = 44511
$ws.Cells.Item(243, 5).Value2 = 8
$ws.Cells.Item(243, 6).Value2 = "Fruta"
$ws.Cells.Item(243, 7).Value2 = 100104
$ws.Cells.Item(243, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(243, 9).Value2 = 100104005
$ws.Cells.Item(243, 10).Value2 = "Pera"
$ws.Cells.Item(243, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(243, 12).Value2 = "Primera"
$ws.Cells.Item(243, 13).Value2 = 240
$ws.Cells.Item(243, 14).Value2 = 10000
$ws.Cells.Item(243, 15).Value2 = 11000
$ws.Cells.Item(243, 16).Value2 = 10417
$ws.Cells.Item(243, 17).Value2 = "`$/caja 17 kilos empedrada"
$ws.Cells.Item(243, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(243, 19).Value2 = 613
$ws.Cells.Item(243, 20).Value2 = 17
